# Reorder the worksheets: move "总计" (the small summary sheet, currently
# second) so that it becomes the first sheet, ahead of "2021-Q1" (the big
# fund-holdings sheet). This matches the commit "update data with resort
# sheetname" / the workbook.xml <sheets> reordering in the diff, where
# "总计" now appears before "2021-Q1" in tab order.

$wb = $excel.ActiveWorkbook

$quarterSheet = $wb.Worksheets.Item("2021-Q1")
$totalSheet   = $wb.Worksheets.Item("总计")

# Move "总计" to be before "2021-Q1" -> new tab order: 总计, 2021-Q1
$totalSheet.Move($quarterSheet)
